$wb = $excel.ActiveWorkbook

# 1. Rename the existing sheet "Parametros" -> "Candidatas"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Candidatas"

# 2. Add a new worksheet right after "Candidatas" and name it "Monitoreadas"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Monitoreadas"

# 3. Populate header row + data on the new sheet
$ws2.Range("B1").Value = "Bus k"
$ws2.Range("C1").Value = "Bus m"
$ws2.Range("D1").Value = "id"
$ws2.Range("E1").Value = "Rating"

$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 101
$ws2.Range("D2").Value = 1
$ws2.Range("E2").Value = 50

$ws2.Range("B3").Value = 1
$ws2.Range("C3").Value = 201
$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = 50

# 4. Copy the header style (bold, centered, thin border) from "Candidatas" B1
$ws1.Range("B1").Copy()
$ws2.Range("B1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 5. Selections on each sheet
[void]$ws1.Range("H1:H3").Select()
[void]$ws2.Range("C3").Select()

# 6. Make "Monitoreadas" the active sheet/tab
[void]$ws2.Activate()
